$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily case data for rows 386-464 (dates 2021-09-21 .. 2021-12-08)
$colA = @(44460,44461,44462,44463,44464,44465,44466,44467,44468,44469,44470,44471,44472,44473,44474,44475,44476,44477,44478,44479,44480,44481,44482,44483,44484,44485,44486,44487,44488,44489,44490,44491,44492,44493,44494,44495,44496,44497,44498,44499,44500,44501,44502,44503,44504,44505,44506,44507,44508,44509,44510,44511,44512,44513,44514,44515,44516,44517,44518,44519,44520,44521,44522,44523,44524,44525,44526,44527,44528,44529,44530,44531,44532,44533,44534,44535,44536,44537,44538)
$colB = @(0,0,0,1,0,0,1,0,0,0,1,0,2,2,3,2,0,0,0,0,0,0,0,0,0,0,3,0,0,0,0,2,0,2,1,0,0,6,1,0,0,1,0,0,0,1,0,0,2,0,1,0,3,0,0,0,5,1,0,1,1,0,2,0,10,4,8,3,1,5,1,0,2,2,3,4,2,4,1)
$colC = @(3,3,2,3,2,1,2,2,2,2,2,2,4,5,8,10,10,9,9,7,5,2,0,0,0,0,3,3,3,3,3,5,5,4,5,5,5,11,10,10,8,8,8,8,2,2,2,2,3,3,4,4,6,6,6,4,9,9,9,7,8,8,10,5,14,18,25,27,28,31,32,22,20,14,14,17,14,17,18)
$colD = @(35.34817956875221,35.34817956875221,23.5654530458348,35.34817956875221,23.5654530458348,11.7827265229174,23.5654530458348,23.5654530458348,23.5654530458348,23.5654530458348,23.5654530458348,23.5654530458348,47.13090609166961,58.91363261458702,94.26181218333922,117.827265229174,117.827265229174,106.0445387062566,106.0445387062566,82.47908566042182,58.91363261458702,23.5654530458348,0,0,0,0,35.34817956875221,35.34817956875221,35.34817956875221,35.34817956875221,35.34817956875221,58.91363261458702,58.91363261458702,47.13090609166961,58.91363261458702,58.91363261458702,58.91363261458702,129.6099917520914,117.827265229174,117.827265229174,94.26181218333922,94.26181218333922,94.26181218333922,94.26181218333922,23.5654530458348,23.5654530458348,23.5654530458348,23.5654530458348,35.34817956875221,35.34817956875221,47.13090609166961,47.13090609166961,70.69635913750442,70.69635913750442,70.69635913750442,47.13090609166961,106.0445387062566,106.0445387062566,106.0445387062566,82.47908566042182,94.26181218333922,94.26181218333922,117.827265229174,58.91363261458702,164.9581713208436,212.0890774125133,294.5681630729351,318.1336161187699,329.9163426416873,365.2645222104395,377.0472487333569,259.2199835041828,235.6545304583481,164.9581713208436,164.9581713208436,200.3063508895958,164.9581713208436,200.3063508895958,212.0890774125133)

$startRow = 386
$endRow = $startRow + $colA.Count - 1

# Extend column A date formatting/style down to the new rows by
# copying the format of the last existing data row (A385).
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $colA.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $colA[$i]
  $ws.Cells.Item($r, 2).Value = $colB[$i]
  $ws.Cells.Item($r, 3).Value = $colC[$i]
  $ws.Cells.Item($r, 4).Value = $colD[$i]
}

$excel.CutCopyMode = 0
